$wb = $excel.ActiveWorkbook
$nl = [char]10

# --- Metadata sheet updates ---
$wsMeta = $wb.Worksheets.Item("Metadata")

# Date property (row 8, col B)
$wsMeta.Range("B8").Value = "2025-06-13T15:45:04+00:00"

# FHIR Version property (row 15, col B)
$wsMeta.Range("B15").Value = "4.0.1"

# --- Elements sheet updates ---
$wsElem = $wb.Worksheets.Item("Elements")

# Row 2 = "Extension" element: Constraint(s) text simplified (drop the
# "unless an empty Parameters resource ... or $this is Parameters" clause)
$wsElem.Range("AJ2").Value = "ele-1:All FHIR elements must have a @value or children {hasValue() or (children().count() > id.count())}" + $nl + "ext-1:Must have either extensions or value[x], not both {extension.exists() != value.exists()}"

# Row 3 = "Extension.id": Type(s) changes from "id" to "string"
$wsElem.Range("K3").Value = "string" + $nl

# Row 6 = "Extension.value[x]": Definition now references R4 (not R4B) docs
$wsElem.Range("M6").Value = "Value of extension - must be one of a constrained set of the data types (see [Extensibility](http://hl7.org/fhir/R4/extensibility.html) for a list)."

# Row 6 = "Extension.value[x]": Binding Value Set now points at the IG's own value set
$wsElem.Range("Z6").Value = "https://nih-ncpi.github.io/ncpi-fhir-ig-2/ValueSet/edam-ontology-terms-vs"

# Column Z widens to fit the new, longer URL value
$wsElem.Columns.Item(26).ColumnWidth = 58.28125
